$wb = $excel.ActiveWorkbook

$wsLevel1 = $wb.Worksheets.Item("level1")
$wsLevel2 = $wb.Worksheets.Item("level2")
$wsLevel3 = $wb.Worksheets.Item("level3")

# level1 sheet changes
$wsLevel1.Range("I27").Value = -180
$wsLevel1.Range("J27").Value = 360

# level2 sheet changes
$wsLevel2.Range("B4").Value = $true
$wsLevel2.Range("F4").Value = 359000
$wsLevel2.Range("H7").Value = 90
$wsLevel2.Range("I27").Value = -180
$wsLevel2.Range("J27").Value = 360

# level3 sheet changes
$wsLevel3.Range("B4").Value = $true
$wsLevel3.Range("F4").Value = 359000
$wsLevel3.Range("F7").Value = 100
$wsLevel3.Range("H7").Value = 90
$wsLevel3.Range("B16").Value = $true
$wsLevel3.Range("F16").Value = "internal gis dataset"
$wsLevel3.Range("B21").Value = $true
$wsLevel3.Range("F21").Value = "internal gis dataset"
$wsLevel3.Range("B26").Value = $true
$wsLevel3.Range("F26").Value = "depends"
$wsLevel3.Range("B27").Value = $true
$wsLevel3.Range("F27").Value = "depends"
$wsLevel3.Range("I27").Value = -180
$wsLevel3.Range("J27").Value = 360
$wsLevel3.Range("B28").Value = $true
$wsLevel3.Range("F28").Value = "depends"

# selections / active sheet
$wsLevel1.Range("I27:J27").Select()
$wsLevel3.Range("F7:J7").Select()
$wsLevel2.Activate()
$wsLevel2.Range("H8").Select()
